$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row for "https://permissionspolicy.xyz/" (row 2); this shifts
# google.com and facebook.com up by one row and drops the now-empty last row.
$ws.Rows(2).Delete()
